$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) store plain text values (inline strings)
# in the source sheet (e.g. "296.98", "-1.56%"). Excel's COM Range.Value
# setter auto-detects numeric-looking / percent-looking strings and would
# silently convert them to real numbers, losing the original textual
# formatting (trailing zeros, exact decimal digits, "%" literal, etc.).
# Marking the range as Text ("@") before writing keeps the values as text,
# and resetting the style back to "Normal" afterwards avoids leaving a
# lingering custom number format on the cells.
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = '296.98'
$ws.Range("E2").Value = '-1.56%'
$ws.Range("D3").Value = '31.70'
$ws.Range("E3").Value = '0.69%'
$ws.Range("D4").Value = '5.065'
$ws.Range("E4").Value = '-1.62%'
$ws.Range("D5").Value = '0.08086'
$ws.Range("E5").Value = '8.98%'
$ws.Range("D6").Value = '2.577'
$ws.Range("E6").Value = '19.93%'
$ws.Range("D7").Value = '7.810'
$ws.Range("E7").Value = '-1.43%'
$ws.Range("D8").Value = '3.822'
$ws.Range("E8").Value = '1.60%'
$ws.Range("D9").Value = '0.9230'
$ws.Range("E9").Value = '-0.33%'
$ws.Range("D10").Value = '0.1756'
$ws.Range("E10").Value = '1.72%'
$ws.Range("D11").Value = '0.07443'
$ws.Range("E11").Value = '-3.18%'
$ws.Range("D12").Value = '0.08981'
$ws.Range("E12").Value = '9.79%'
$ws.Range("D13").Value = '0.03039'
$ws.Range("E13").Value = '0.40%'
$ws.Range("E14").Value = '0.79%'
$ws.Range("D15").Value = '0.001498'
$ws.Range("E15").Value = '0.62%'
$ws.Range("D16").Value = '0.006016'
$ws.Range("E16").Value = '-1.82%'
$ws.Range("D17").Value = '3.553'
$ws.Range("D18").Value = '2.247'
$ws.Range("E18").Value = '0.71%'
$ws.Range("E19").Value = '0.29%'
$ws.Range("E20").Value = '0.07%'
$ws.Range("D21").Value = '4.008'
$ws.Range("E21").Value = '-13.77%'
$ws.Range("D23").Value = '0.04594'
$ws.Range("E23").Value = '-1.01%'
$ws.Range("D24").Value = '0.001241'
$ws.Range("E24").Value = '1.64%'
$ws.Range("E25").Value = '-1.09%'
$ws.Range("E26").Value = '-7.66%'
$ws.Range("D27").Value = '0.0003406'
$ws.Range("E27").Value = '82.00%'
$ws.Range("D39").Value = '0.01772'
$ws.Range("E39").Value = '1.83%'
$ws.Range("D40").Value = '0.04505'
$ws.Range("E40").Value = '-0.99%'
$ws.Range("D41").Value = '0.006842'
$ws.Range("E41").Value = '-4.02%'
$ws.Range("D42").Value = '0.1351'
$ws.Range("E42").Value = '0.16%'
$ws.Range("D43").Value = '0.002207'
$ws.Range("E43").Value = '0.03%'
$ws.Range("D44").Value = '0.009835'
$ws.Range("E44").Value = '-10.19%'
$ws.Range("D45").Value = '0.00006456'
$ws.Range("E45").Value = '2.87%'
$ws.Range("E46").Value = '-0.19%'
$ws.Range("D47").Value = '0.008734'
$ws.Range("E47").Value = '24.82%'
$ws.Range("E48").Value = '11.15%'
$ws.Range("E49").Value = '-0.19%'
$ws.Range("D50").Value = '0.0001996'
$ws.Range("E50").Value = '-0.12%'

$priceVolumeRange.Style = "Normal"
